$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6639
$ws.Range("F3").Value = 786
$ws.Range("F4").Value = 1100
$ws.Range("F6").Value = 667
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = 768
$ws.Range("F11").Value = 706
$ws.Range("F12").Value = 1296
$ws.Range("F14").Value = 105
$ws.Range("F15").Value = 522
$ws.Range("F16").Value = 524
$ws.Range("F19").Value = 1049
$ws.Range("F20").Value = 1468
$ws.Range("F22").Value = 436
$ws.Range("F23").Value = 434
$ws.Range("F25").Value = 1111
$ws.Range("F26").Value = 230
$ws.Range("F27").Value = 2329
$ws.Range("F28").Value = 260
$ws.Range("F29").Value = 881
$ws.Range("F30").Value = 429
$ws.Range("F32").Value = 3776

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 36
$ws.Range("F6").Value = 741
$ws.Range("F11").Value = 145
$ws.Range("F12").Value = 639
$ws.Range("F13").Value = 4
$ws.Range("F17").Value = 391
$ws.Range("F18").Value = 325
$ws.Range("F19").Value = 4114
$ws.Range("F24").Value = 224
$ws.Range("F26").Value = 102
$ws.Range("F28").Value = 221
$ws.Range("F29").Value = 39
$ws.Range("F32").Value = 1697

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 92
$ws.Range("F4").Value = 1233
$ws.Range("F5").Value = 1615
$ws.Range("F8").Value = 928

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 92
$ws.Range("F3").Value = 1233
$ws.Range("F4").Value = 1615
$ws.Range("F7").Value = 928
$ws.Range("F9").Value = 6639
$ws.Range("F10").Value = 36
$ws.Range("F11").Value = 786
$ws.Range("F12").Value = 741
$ws.Range("F14").Value = 667
$ws.Range("F15").Value = 667
$ws.Range("F17").Value = 768
$ws.Range("F21").Value = 145
$ws.Range("F22").Value = 145
$ws.Range("F23").Value = 4
$ws.Range("F24").Value = 1296
$ws.Range("F26").Value = 105
$ws.Range("F27").Value = 522
$ws.Range("F28").Value = 524
$ws.Range("F30").Value = 391
$ws.Range("F32").Value = 1049
$ws.Range("F33").Value = 1468
$ws.Range("F36").Value = 436
$ws.Range("F37").Value = 434
$ws.Range("F40").Value = 224
$ws.Range("F41").Value = 1111
$ws.Range("F42").Value = 230
$ws.Range("F43").Value = 2329
$ws.Range("F44").Value = 39
$ws.Range("F45").Value = 1697
$ws.Range("F46").Value = 1697
$ws.Range("F47").Value = 881
$ws.Range("F48").Value = 429
$ws.Range("F49").Value = 3776
